$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 18551
$ws.Range("I34").Value = 18551
$ws.Range("K34").Value = 18551
$ws.Range("M34").Value = -18348
$ws.Range("H36").Value = 18551
$ws.Range("I36").Value = 18551
$ws.Range("K36").Value = 18551
$ws.Range("M36").Value = -17836
$ws.Range("H53").Value = 94.125
$ws.Range("I53").Value = 76
$ws.Range("J53").Value = 100.166664
$ws.Range("K53").Value = 76
$ws.Range("L53").Value = 100.166664
$ws.Range("M53").Value = 561
$ws.Range("N53").Value = -1374.166664
$ws.Range("H107").Value = 656.8889
$ws.Range("I107").Value = 433.5
$ws.Range("J107").Value = 768.5833
$ws.Range("K107").Value = 433.5
$ws.Range("L107").Value = 768.5833
$ws.Range("M107").Value = 1486.5
$ws.Range("N107").Value = -4608.5833
$ws.Range("H111").Value = 2312.2222
$ws.Range("I111").Value = 2370.3333
$ws.Range("J111").Value = 2283.1667
$ws.Range("K111").Value = 7110.999899999999
$ws.Range("L111").Value = 6849.500100000001
$ws.Range("M111").Value = -4043.999899999999
$ws.Range("N111").Value = -12983.5001
$ws.Range("H113").Value = 2550.5789
$ws.Range("J113").Value = 3540.8572
$ws.Range("L113").Value = 3540.8572
$ws.Range("N113").Value = -10048.8572
$ws.Range("H137").Value = 1266.9286
$ws.Range("I137").Value = 1104.7894
$ws.Range("J137").Value = 1609.2222
$ws.Range("K137").Value = 3314.3682
$ws.Range("L137").Value = 4827.6666
$ws.Range("M137").Value = -764.3681999999999
$ws.Range("N137").Value = -9927.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9497.397000000001
$ws.Range("I32").Value = 9087.172
$ws.Range("J32").Value = 16061
$ws.Range("K32").Value = 9087.172
$ws.Range("L32").Value = 16061
$ws.Range("M32").Value = -8800.172
$ws.Range("N32").Value = -16635
$ws.Range("H45").Value = 1124.875
$ws.Range("I45").Value = 1166.6666
$ws.Range("J45").Value = 999.5
$ws.Range("K45").Value = 1166.6666
$ws.Range("L45").Value = 999.5
$ws.Range("M45").Value = -789.6666
$ws.Range("N45").Value = -1753.5
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 29852.059
$ws.Range("J62").Value = 30150.605
$ws.Range("L62").Value = 30150.605
$ws.Range("N62").Value = -31522.605
$ws.Range("H65").Value = 29852.059
$ws.Range("J65").Value = 30150.605
$ws.Range("L65").Value = 90451.815
$ws.Range("N65").Value = -97315.815
$ws.Range("H102").Value = 178333.33
$ws.Range("I102").Value = 178333.33
$ws.Range("K102").Value = 178333.33
$ws.Range("M102").Value = -175088.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2460.0857
$ws.Range("I31").Value = 1400.1428
$ws.Range("J31").Value = 4050
$ws.Range("K31").Value = 1400.1428
$ws.Range("L31").Value = 4050
$ws.Range("M31").Value = -1105.1428
$ws.Range("N31").Value = -4640
$ws.Range("H34").Value = 2460.0857
$ws.Range("I34").Value = 1400.1428
$ws.Range("J34").Value = 4050
$ws.Range("K34").Value = 1400.1428
$ws.Range("L34").Value = 4050
$ws.Range("M34").Value = -1198.1428
$ws.Range("N34").Value = -4454
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H107").Value = 345.42856
$ws.Range("I107").Value = 299.05884
$ws.Range("J107").Value = 417.0909
$ws.Range("K107").Value = 299.05884
$ws.Range("L107").Value = 417.0909
$ws.Range("M107").Value = 1620.94116
$ws.Range("N107").Value = -4257.0909

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 526.2174
$ws.Range("J113").Value = 507.46155
$ws.Range("L113").Value = 1522.38465
$ws.Range("N113").Value = -5862.38465
$ws.Range("H118").Value = 2704.2
$ws.Range("I118").Value = 1283.3334
$ws.Range("J118").Value = 3313.1428
$ws.Range("K118").Value = 3850.0002
$ws.Range("L118").Value = 9939.428400000001
$ws.Range("M118").Value = -2607.0002
$ws.Range("N118").Value = -12425.4284
$ws.Range("H131").Value = 4391.3657
$ws.Range("I131").Value = 11534.444
$ws.Range("J131").Value = 2382.375
$ws.Range("K131").Value = 34603.33199999999
$ws.Range("L131").Value = 7147.125
$ws.Range("M131").Value = -29563.33199999999
$ws.Range("N131").Value = -17227.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2548.0908
$ws.Range("I122").Value = 2022.2222
$ws.Range("J122").Value = 2912.1538
$ws.Range("K122").Value = 6066.6666
$ws.Range("L122").Value = 8736.4614
$ws.Range("M122").Value = -3616.6666
$ws.Range("N122").Value = -13636.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1043.75
$ws.Range("I22").Value = 941.6667
$ws.Range("J22").Value = 1350
$ws.Range("K22").Value = 941.6667
$ws.Range("L22").Value = 1350
$ws.Range("M22").Value = -646.6667
$ws.Range("N22").Value = -1940
$ws.Range("H27").Value = 1043.75
$ws.Range("I27").Value = 941.6667
$ws.Range("J27").Value = 1350
$ws.Range("K27").Value = 941.6667
$ws.Range("L27").Value = 1350
$ws.Range("M27").Value = -834.6667
$ws.Range("N27").Value = -1564
$ws.Range("H46").Value = 1483.3334
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 1625
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 1625
$ws.Range("M46").Value = -1012
$ws.Range("N46").Value = -2001
$ws.Range("H122").Value = 25003664
$ws.Range("I122").Value = 17860638
$ws.Range("J122").Value = 50004252
$ws.Range("K122").Value = 53581914
$ws.Range("L122").Value = 150012756
$ws.Range("M122").Value = -53579464
$ws.Range("N122").Value = -150017656
$ws.Range("H125").Value = 24475
$ws.Range("J125").Value = 24475
$ws.Range("L125").Value = 24475
$ws.Range("N125").Value = -34315
$ws.Range("H132").Value = 2892.8057
$ws.Range("I132").Value = 2076.5356
$ws.Range("K132").Value = 6229.6068
$ws.Range("M132").Value = -3699.6068
$ws.Range("H136").Value = 16834994
$ws.Range("I136").Value = 21740402
$ws.Range("J136").Value = 717222.5
$ws.Range("K136").Value = 65221206
$ws.Range("L136").Value = 2151667.5
$ws.Range("M136").Value = -65218656
$ws.Range("N136").Value = -2156767.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4200
$ws.Range("H65").Value = 4200
$ws.Range("H96").Value = 2031.9286
$ws.Range("I96").Value = 1491.1666
$ws.Range("J96").Value = 2437.5
$ws.Range("K96").Value = 1491.1666
$ws.Range("L96").Value = 2437.5
$ws.Range("M96").Value = -118.1666
$ws.Range("N96").Value = -5183.5
$ws.Range("H109").Value = 40188
$ws.Range("J109").Value = 40188
$ws.Range("L109").Value = 40188
$ws.Range("N109").Value = -42962
$ws.Range("H122").Value = 8225757
$ws.Range("I122").Value = 9616713
$ws.Range("J122").Value = 5212018
$ws.Range("K122").Value = 28850139
$ws.Range("L122").Value = 15636054
$ws.Range("M122").Value = -28847689
$ws.Range("N122").Value = -15640954
$ws.Range("H132").Value = 1736.6666
$ws.Range("I132").Value = 1097.5714
$ws.Range("J132").Value = 3363.4546
$ws.Range("K132").Value = 3292.7142
$ws.Range("L132").Value = 10090.3638
$ws.Range("M132").Value = -762.7142000000003
$ws.Range("N132").Value = -15150.3638
